$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 30 contains the "Total" row, which should be removed.
# Deleting the entire row shifts rows 31-33 (Xinjiang, Yunnan, Zhejiang) up
# by one, leaving the correct data in place and dropping the now-duplicate
# last row automatically.
$ws.Rows.Item(30).Delete()
